# update to manual status column
# The manualStatus column (I) for the first three data rows (5-7) used to
# hold the raw numeric index count (4); it is switched to the textual
# "[4]" tag used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "[4]"
$ws.Range("I6").Value = "[4]"
$ws.Range("I7").Value = "[4]"

# Give the fastqFileName column (F) enough room to show the full file name.
$ws.Columns.Item(6).ColumnWidth = 55.16

# Rows 6/7 shrink slightly to their content height after the edit.
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# Leave the active selection on the last touched cell.
$ws.Range("I7").Select()
